$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.649903655052185
$ws.Range("B1").Value = 3.513938426971436
$ws.Range("C1").Value = 4.098756313323975
$ws.Range("D1").Value = 1.288067936897278
$ws.Range("E1").Value = 0.7553884983062744
